# Aggiornamento del glossario: nuovi id statici (ITCH...) e nuova data di
# inserimento entry (serial 43538) per le 50 righe dati, con relativo
# aggiornamento dei formati data (numFmt 164/165) per includere l'orario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nuovo formato data/ora per le celle della colonna "Data_inserimento_entry"
$ws.Range("B2:B51").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$newDateSerial = 43538
$idPrefix = "ITCH158255732172"
$idStart = 59687

for ($i = 0; $i -lt 50; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDateSerial
    $ws.Cells.Item($row, 3).Value = $idPrefix + ($idStart + $i)
}
